$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q1)
$ws.Range("C2").Value = 22.96
$ws.Range("E2").Value = 18

# Row 3 (Q3)
$ws.Range("C3").Value = 6.49
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 5
$ws.Range("J3").Value = 11.94
$ws.Range("L3").Value = "{'Q1': np.float64(4.0), 'Q2': np.float64(6.0), 'Q3': np.float64(10.0), 'Q4': np.float64(12.0)}"

# Row 4 (Q4)
$ws.Range("C4").Value = 4.45
$ws.Range("I4").Value = 2.28
$ws.Range("J4").Value = 5.21

# Row 5 (Q7)
$ws.Range("C5").Value = 4.01
$ws.Range("E5").Value = 6
$ws.Range("I5").Value = 2.59
$ws.Range("J5").Value = 6.73

# Row 6 (Q9)
$ws.Range("C6").Value = 14.09
$ws.Range("E6").Value = 23
$ws.Range("I6").Value = 6.7
$ws.Range("J6").Value = 44.93

# Row 7 (Q11)
$ws.Range("E7").Value = 1
$ws.Range("J7").Value = 6.6

# Row 8 (Q13)
$ws.Range("C8").Value = 5.52
$ws.Range("E8").Value = 5
$ws.Range("I8").Value = 2.87
$ws.Range("J8").Value = 8.26
$ws.Range("K8").Value = 0.52
